$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gens")
$ws.Range("C4").Value = 76
$ws.Range("C5").Value = 37.406531
$ws.Range("C8").Value = 76
$ws.Range("C20").Value = 12
$ws.Range("C21").Value = 12
$ws = $wb.Worksheets.Item("lines")
$ws.Range("C2").Value = 39.880351
$ws.Range("D2").Value = -0.025674804
$ws.Range("C3").Value = -55.607927
$ws.Range("D3").Value = -0.64420583
$ws.Range("C4").Value = 21.134107
$ws.Range("D4").Value = 0.41539694
$ws.Range("C5").Value = -1.4590775
$ws.Range("D5").Value = -0.67019363
$ws.Range("C6").Value = 20.339429
$ws.Range("D6").Value = 0.66109451
$ws.Range("C7").Value = 35.78614
$ws.Range("D7").Value = -0.60048258
$ws.Range("C8").Value = -271.39407
$ws.Range("D8").Value = 0.16740896
$ws.Range("C9").Value = -75.45907800000001
$ws.Range("D9").Value = -0.54881998
$ws.Range("C10").Value = -49.865893
$ws.Range("D10").Value = 0.43005801
$ws.Range("C11").Value = -115.66057
$ws.Range("D11").Value = 0.21003524
$ws.Range("E11").Value = 0
$ws.Range("C12").Value = -149.40653
$ws.Range("F12").Value = 0
$ws.Range("C13").Value = -175
$ws.Range("D13").Value = 5.6445451
$ws.Range("E13").Value = 13.379234
$ws.Range("C14").Value = -145.40653
$ws.Range("D14").Value = -5.6445451
$ws.Range("C15").Value = -178.44467
$ws.Range("D15").Value = 0.99847946
$ws.Range("C16").Value = -211.22826
$ws.Range("D16").Value = 1.0079594
$ws.Range("C17").Value = -236.5747
$ws.Range("D17").Value = -1.0916639
$ws.Range("C18").Value = -269.35829
$ws.Range("D18").Value = -1.0821839
$ws.Range("C19").Value = -276.76813
$ws.Range("D19").Value = -0.0060025472
$ws.Range("C20").Value = -138.25125
$ws.Range("D20").Value = -0.041339985
$ws.Range("C21").Value = -219.39685
$ws.Range("D21").Value = -0.015482493
$ws.Range("C22").Value = -261.18971
$ws.Range("D22").Value = -0.054424128
$ws.Range("C23").Value = -170.16498
$ws.Range("D23").Value = -0.038941635
$ws.Range("C24").Value = -332.25125
$ws.Range("D24").Value = -0.058072836
$ws.Range("C25").Value = 45.259712
$ws.Range("D25").Value = 0.027540981
$ws.Range("C26").Value = -227.32689
$ws.Range("D26").Value = 0.009136199899999999
$ws.Range("C27").Value = -227.32689
$ws.Range("D27").Value = 0.009136199899999999
$ws.Range("C28").Value = 271.39407
$ws.Range("D28").Value = -0.10363412
$ws.Range("C29").Value = -312.34622
$ws.Range("D29").Value = -0.009695559100000001
$ws.Range("E29").Value = 0
$ws.Range("C30").Value = 80.35468400000001
$ws.Range("D30").Value = 0.023199605
$ws.Range("C31").Value = -172.53448
$ws.Range("D31").Value = -0.0045158931
$ws.Range("C32").Value = -139.81174
$ws.Range("D32").Value = -0.0052859442
$ws.Range("C33").Value = -52.76724
$ws.Range("D33").Value = -0.0041933293
$ws.Range("C34").Value = -52.76724
$ws.Range("D34").Value = -0.0041933293
$ws.Range("C35").Value = -50.322658
$ws.Range("D35").Value = 0.02017357
$ws.Range("C36").Value = -50.322658
$ws.Range("D36").Value = 0.02017357
$ws.Range("C37").Value = -114.32266
$ws.Range("D37").Value = 0.011095463
$ws.Range("C38").Value = -114.32266
$ws.Range("D38").Value = 0.011095463
$ws.Range("C39").Value = -160.18826
$ws.Range("D39").Value = 0.0034232782
$ws = $wb.Worksheets.Item("bus")
$ws.Range("B2").Value = 108.6
$ws.Range("B3").Value = 108.57433
$ws.Range("C3").Value = -0.55832492
$ws.Range("B4").Value = 107.95579
$ws.Range("C4").Value = 11.733273
$ws.Range("B5").Value = 107.90413
$ws.Range("C5").Value = -0.37302208
$ws.Range("B6").Value = 109.0154
$ws.Range("C6").Value = -1.7963991
$ws.Range("B7").Value = 109.23542
$ws.Range("C7").Value = -4.4634953
$ws.Range("B8").Value = 115.09
$ws.Range("C8").Value = -30.514076
$ws.Range("B9").Value = 115.09
$ws.Range("C9").Value = -21.400278
$ws.Range("B10").Value = 107.35531
$ws.Range("C10").Value = 7.474722
$ws.Range("B11").Value = 109.44545
$ws.Range("C11").Value = 2.5917996
$ws.Range("B12").Value = 108.35379
$ws.Range("C12").Value = 22.464075
$ws.Range("B13").Value = 108.36327
$ws.Range("C13").Value = 25.217896
$ws.Range("B14").Value = 108.34779
$ws.Range("C14").Value = 35.748945
$ws.Range("B15").Value = 108.31245
$ws.Range("C15").Value = 28.270627
$ws.Range("B16").Value = 108.22684
$ws.Range("C16").Value = 48.642866
$ws.Range("B17").Value = 108.25438
$ws.Range("C17").Value = 47.873451
$ws.Range("B18").Value = 108.24468
$ws.Range("C18").Value = 55.994452
$ws.Range("B19").Value = 108.24017
$ws.Range("C19").Value = 58.409935
$ws.Range("B20").Value = 108.27758
$ws.Range("C20").Value = 46.025293
$ws.Range("B21").Value = 108.29775
$ws.Range("C21").Value = 48.038199
$ws.Range("B22").Value = 108.23597
$ws.Range("C22").Value = 59.781883
$ws.Range("B23").Value = 108.2394
$ws.Range("C23").Value = 70.674685
$ws.Range("B24").Value = 108.30885
$ws.Range("C24").Value = 50.553298
$ws.Range("B25").Value = 108.1232
$ws.Range("C25").Value = 34.530374
